# Update "想去人数" (number of people wanting to go) figures that changed
# between scrapes on the 展览 (Exhibition) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F5").Value = 829
    $ws.Range("F6").Value = 8
    $ws.Range("F7").Value = 420
}
